# Add a new diary entry (row 32) to the "Journal" worksheet:
#   Date: 2024-08-24 (serial 45528), Start: 08:30, End: 14:28,
#   Category: "Réalisation de l'application " (existing shared string)
# The Duration (column D) cell already contains a shared formula that
# will recompute automatically once B32/C32 are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

$ws.Range("A32").Value = 45528
$ws.Range("B32").Value = 0.35416666666666669
$ws.Range("C32").Value = 0.60277777777777775
$ws.Range("E32").Value = "Réalisation de l'application "

$excel.Calculate()

# Reflect the author's last active selection when the file was saved.
$ws.Range("H26").Select() | Out-Null
